$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Luettavuus" (readability) column entirely. This shifts the
#    old D..H columns left into C..G and drops the shared strings that only
#    lived in that column.
$ws.Columns("C").Delete()

# 2. Re-apply the column widths for the two columns whose content changed
#    shape now that "Sisalto lyhyesti" (ex-D, now C) and "Hyvaa" (ex-E, now D)
#    are adjacent to the removed column.
$ws.Columns("C").ColumnWidth = 54.3
$ws.Columns("D").ColumnWidth = 33.7

# 3. Add the two new article rows (spacer blank rows at 11 and 13 keep the
#    same every-other-row layout the sheet already used for rows 5/7/9).
$ws.Range("A12").Value = "Instructional Design of a Programming Course -- A Learning Theoretic Approach"
$ws.Range("B12").Value = "Käsittelee ohjelmoinnin opettamista muun muossa CA:ta käyttäen Aarhusin yliopistolla 5/5"
$ws.Range("C12").Value = "Kertoo kolmesta eri opetustekniikasta ohjelmoinnissa (cognitive load theory, CA, worked examples). Selitetään miten kutakin voi hyödyntää ohjelmoinnin opettamiseen ja kerrotaan miten Aarhusin yliopistolla näitä on hyödynnetty opetuksessa."
$ws.Range("D12").Value = "Käy todella hyvin erilaisia opetustekniikoita läpi ja selitetään miten niitä konkreettisesti opetuksessa ja materiaaleja suunniteltaessa voidaan hyödyntää"
$ws.Range("E12").Value = "Guidance-fading negatiivinen vaikutus henkilöille, jotka jo omaksuneet expertisen. "
$ws.Range("F12").Value = "400 oppilasta per vuosi, Aarhuusin yliopisto, Tanska, opetusmenetelmät olleet käytössä yli 4 vuotta kyseisessä yliopistossa"
$ws.Range("G12").Value = "Ei ole kunnon tuloksia, mutta toiminut hyvin. Pitäisi tehdä kontrolloituja testejä tulosten muodostamiseen."

$ws.Range("A14").Value = "Helping Novice Programming Students Succeed"
$ws.Range("B14").Value = "Käsittelee ohjelmoinnin opettamista käyttäen CA:ta, University of New Mexico, Valencia Campus 5/5"
$ws.Range("C14").Value = "Kertoo aluksi pääpiirteittäin mitä on CA ja tämän jälkeen syventyy siihen, miten CA:ta voidaan hyödyntää ohjelmoinnin opettamiseen"
$ws.Range("D14").Value = "Käy läpi oleellista asiaa, CA:n pääpiirteittäin ja miten sitä voi hyödyntää ohjelmoinnin opettamisessa"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "2006, University of New Mexico, Valencia Campus"
$ws.Range("G14").Value = "Ei tuloksia, aikoivat juuri ottaa käyttöön CA:n ohjelmoinnin opetuksessa"
